$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.270.80'
$ws.Range('D3').Value = '2.175.89'
$ws.Range('E3').Value = '  -1.37%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '''237.91'
$ws.Range('E5').Value = '  -1.06%  '
$ws.Range('E6').Value = '  -1.70%  '
$ws.Range('D7').Value = '''70.32'
$ws.Range('E7').Value = '  -3.41%  '
$ws.Range('D9').Value = '''0.581'
$ws.Range('E9').Value = '  -3.80%  '
$ws.Range('D10').Value = '''40.36'
$ws.Range('E10').Value = '  -5.65%  '
$ws.Range('D11').Value = '''0.0927'
$ws.Range('E11').Value = '  -2.12%  '
$ws.Range('D12').Value = '''54.34'
$ws.Range('E12').Value = '  -5.52%  '
$ws.Range('D13').Value = '''6.78'
$ws.Range('E13').Value = '  -4.19%  '
$ws.Range('E14').Value = '  -1.93%  '
$ws.Range('D15').Value = '2.499.06'
$ws.Range('E15').Value = '  -1.53%  '
$ws.Range('D16').Value = '''13.99'
$ws.Range('E16').Value = '  -1.10%  '
$ws.Range('E17').Value = '  -4.02%  '
$ws.Range('D18').Value = '2.178.22'
$ws.Range('E18').Value = '  -1.51%  '
$ws.Range('D19').Value = '41.090.04'
$ws.Range('E19').Value = '  -1.75%  '
$ws.Range('D20').Value = '''0.0000101'
$ws.Range('E20').Value = '  -6.50%  '
$ws.Range('D21').Value = '''70.56'
$ws.Range('E21').Value = '  -2.69%  '
$ws.Range('E22').Value = '  -2.35%  '
$ws.Range('D23').Value = '''9.87'
$ws.Range('E23').Value = '  -3.07%  '
$ws.Range('D24').Value = '''226.33'
$ws.Range('E24').Value = '  -0.96%  '
$ws.Range('E25').Value = '  -5.41%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').Value = '''10.92'
$ws.Range('E27').Value = '  -4.33%  '
$ws.Range('E28').Value = '  -1.31%  '
$ws.Range('E29').Value = '  -2.10%  '
$ws.Range('E30').Value = '  +0.84%  '
$ws.Range('D31').Value = '''167.68'
$ws.Range('E31').Value = '  +0.28%  '
$ws.Range('D32').Value = '''19.97'
$ws.Range('E32').Value = '  -2.49%  '
$ws.Range('D33').Value = '''31.32'
$ws.Range('E33').Value = '  +8.59%  '
$ws.Range('D34').Value = '''0.0769'
$ws.Range('E34').Value = '  -2.24%  '
$ws.Range('E35').Value = '  -6.57%  '
$ws.Range('E36').Value = '  -2.80%  '
$ws.Range('E37').Value = '  -6.38%  '
$ws.Range('E38').Value = '  -2.52%  '
$ws.Range('E39').Value = '  -4.35%  '
$ws.Range('B40').Value = 'LidoDAOToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D40').Value = '''2.09'
$ws.Range('E40').Value = '  -1.25%  '
$ws.Range('B41').Value = 'Celestia'
$ws.Range('C41').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D41').Value = '''11.91'
$ws.Range('E41').Value = '  -7.03%  '
$ws.Range('E42').Value = '  -2.81%  '
$ws.Range('D43').Value = '''60.08'
$ws.Range('E43').Value = '  -7.88%  '
$ws.Range('E44').Value = '  -2.96%  '
$ws.Range('E45').Value = '  -2.38%  '
$ws.Range('E46').Value = '  -4.32%  '
$ws.Range('D47').Value = '''98.43'
$ws.Range('E47').Value = '  -4.92%  '
$ws.Range('E48').Value = '  -1.61%  '
$ws.Range('E49').Value = '  -1.62%  '
$ws.Range('E50').Value = '  -6.91%  '
$ws.Range('E51').Value = '  -2.77%  '
